$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/27/2023  Through  3/5/2023"

# --- Plain numeric updates (style/format unchanged) ---
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = -75
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -85.714285714285
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 4
$ws.Range("L15").Value = 300
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = -71.428571428571
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -40.740740740740
$ws.Range("I16").Value = 31
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -36.734693877551
$ws.Range("L16").Value = 82.352941176470
$ws.Range("M16").Value = -29.545454545454
$ws.Range("N16").Value = -80.625
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -5.714285714285
$ws.Range("I17").Value = 76
$ws.Range("J17").Value = 75
$ws.Range("K17").Value = 1.333333333333
$ws.Range("L17").Value = 15.151515151515
$ws.Range("M17").Value = 68.888888888888
$ws.Range("N17").Value = -44.117647058823
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -48.275862068965
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = -25.490196078431
$ws.Range("L18").Value = 2.702702702702
$ws.Range("M18").Value = 65.217391304347
$ws.Range("N18").Value = -77.906976744186
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 51
$ws.Range("J19").Value = 61
$ws.Range("K19").Value = -16.393442622950
$ws.Range("L19").Value = 6.25
$ws.Range("M19").Value = 13.333333333333
$ws.Range("N19").Value = -31.081081081081
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 21
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -80.645161290322
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -5.714285714285
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = -24.817518248175
$ws.Range("I21").Value = 213
$ws.Range("J21").Value = 265
$ws.Range("K21").Value = -19.622641509434
$ws.Range("L21").Value = 18.333333333333
$ws.Range("M21").Value = 24.561403508771
$ws.Range("N21").Value = -65.92
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -25
$ws.Range("M22").Value = 50
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = -22.727272727272
$ws.Range("I23").Value = 39
$ws.Range("J23").Value = 39
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 34.482758620689
$ws.Range("M23").Value = 44.444444444444
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 41.666666666666
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = -5.797101449275
$ws.Range("I24").Value = 153
$ws.Range("J24").Value = 151
$ws.Range("K24").Value = 1.324503311258
$ws.Range("L24").Value = 19.53125
$ws.Range("M24").Value = 42.990654205607
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 7.692307692307
$ws.Range("I25").Value = 104
$ws.Range("J25").Value = 87
$ws.Range("K25").Value = 19.540229885057
$ws.Range("L25").Value = 42.465753424657
$ws.Range("M25").Value = -26.241134751773
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 300
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = -57.142857142857
$ws.Range("L27").Value = -62.5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 25
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 25
$ws.Range("N28").Value = -75
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = 33.333333333333
$ws.Range("M29").Value = 33.333333333333
$ws.Range("N29").Value = -76.470588235294

# --- Cells that change between numeric and text representation ---
# (copy number format from a stable donor cell of the target style, then set the value;
#  numeric-looking text values are entered with a leading apostrophe so Excel keeps them as text)
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value = 1

$ws.Range("C15").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Value = 1

$ws.Range("C15").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = 1

$ws.Range("K15").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = 0

$ws.Range("C15").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value = 1

$ws.Range("C15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 1

$ws.Range("K15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = 0

$ws.Range("D14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = "'0"

$ws.Range("C15").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 1

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = "'0"

$ws.Range("D14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "***.*"

$ws.Range("C15").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Value = 2

$ws.Range("C15").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Value = 1

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = "'0"

$ws.Range("D14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "***.*"

$ws.Range("C15").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = 2

